# Scheduled market-data refresh: push updated currentAveragePrice(NQ/HQ)
# and LevePrice/LeveProfit figures into each job sheet's leve table.
# Columns: H=currentAveragePrice I=currentAveragePriceNQ J=currentAveragePriceHQ
#          K=LevePriceNQ L=LevePriceHQ M=LeveProfitNQ N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 116: Growing Up
$ws.Range("H116").Value = 7651
$ws.Range("I116").Value = 4559.4
$ws.Range("K116").Value = 4559.4
$ws.Range("M116").Value = -1117.4

# Row 138: All-night Crafting
$ws.Range("H138").Value = 3547.3333

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 2665.3333
$ws.Range("I2").Value = 2499.5
$ws.Range("K2").Value = 2499.5
$ws.Range("M2").Value = -2386.5

# Row 23: A Well-rounded Crew
$ws.Range("H23").Value = 1000
$ws.Range("J23").Value = 1000
$ws.Range("L23").Value = 1000
$ws.Range("N23").Value = -1518

# Row 116: No Scope
$ws.Range("H116").Value = 2665.3333
$ws.Range("I116").Value = 2499.5
$ws.Range("K116").Value = 2499.5
$ws.Range("M116").Value = -205.5

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 986.5714
$ws.Range("I122").Value = 1047.5385
$ws.Range("J122").Value = 194
$ws.Range("K122").Value = 3142.6155
$ws.Range("L122").Value = 582
$ws.Range("M122").Value = -692.6155000000003
$ws.Range("N122").Value = -5482

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 2665.3333
$ws.Range("I3").Value = 2499.5
$ws.Range("K3").Value = 2499.5
$ws.Range("M3").Value = -2385.5

# Row 7: Thank You for Your Business
$ws.Range("H7").Value = 28495
$ws.Range("J7").Value = 28495
$ws.Range("L7").Value = 28495
$ws.Range("N7").Value = -28721

# Row 23: Get a Little Bit Closer
$ws.Range("H23").Value = 500
$ws.Range("I23").Value = 500
$ws.Range("K23").Value = 500
$ws.Range("M23").Value = -217

# Row 33: Mors Dagger
$ws.Range("H33").Value = 11331.6
$ws.Range("I33").Value = 11331.6
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 11331.6
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -10995.6

# Row 55: Streamlining Operations
$ws.Range("H55").Value = 100000
$ws.Range("J55").Value = 100000
$ws.Range("L55").Value = 100000
$ws.Range("N55").Value = -100546

$ws = $wb.Worksheets.Item("CRP")
# Row 13: Compulsory Conjury
$ws.Range("H13").Value = 6832.6665
$ws.Range("I13").Value = 500
$ws.Range("K13").Value = 500
$ws.Range("M13").Value = -361

# Row 19: Shielding Sales
$ws.Range("H19").Value = 13400225
$ws.Range("I19").Value = 14357348
$ws.Range("K19").Value = 14357348
$ws.Range("M19").Value = -14357178

# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 5958.32
$ws.Range("J22").Value = 9383.571
$ws.Range("L22").Value = 9383.571
$ws.Range("N22").Value = -10083.571

# Row 24: What You Need
$ws.Range("H24").Value = 13400225
$ws.Range("I24").Value = 14357348
$ws.Range("K24").Value = 14357348
$ws.Range("M24").Value = -14357178

# Row 31: Wall Not Found
$ws.Range("H31").Value = 3242.111
$ws.Range("I31").Value = 3022.375
$ws.Range("K31").Value = 3022.375
$ws.Range("M31").Value = -2727.375

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 3242.111
$ws.Range("I34").Value = 3022.375
$ws.Range("K34").Value = 3022.375
$ws.Range("M34").Value = -2820.375

# Row 37: Heal Away
$ws.Range("H37").Value = 13744.5
$ws.Range("I37").Value = 13744.5
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 13744.5
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -13637.5

# Row 45: A Tree Grew in Gridania
$ws.Range("H45").Value = 50000
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

# Row 74: License to Heal
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("N74").Value = 0

# Row 77: Purified Polyrhythm (L)
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("N77").Value = 0

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water
$ws.Range("H4").Value = 2311696.5
$ws.Range("I4").Value = 3695517.8
$ws.Range("J4").Value = 5328.1665
$ws.Range("K4").Value = 11086553.4
$ws.Range("L4").Value = 15984.4995
$ws.Range("M4").Value = -11086441.4
$ws.Range("N4").Value = -16208.4995

# Row 58: Bread in the Clouds
$ws.Range("H58").Value = 1531
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

# Row 104: Fits to a Tea
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").ClearContents()
$ws.Range("N104").Value = 0

# Row 107: Slippery Service
$ws.Range("H107").Value = 571.5714
$ws.Range("I107").Value = 483
$ws.Range("J107").Value = 638
$ws.Range("K107").Value = 1449
$ws.Range("L107").Value = 1914
$ws.Range("M107").Value = 471
$ws.Range("N107").Value = -5754

# Row 122: Salt of the North
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 36: Keep the Change
$ws.Range("H36").Value = 999
$ws.Range("I36").Value = 999
$ws.Range("K36").Value = 999
$ws.Range("M36").Value = -514

# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 4533.3335
$ws.Range("I102").Value = 4533.3335
$ws.Range("K102").Value = 4533.3335
$ws.Range("M102").Value = -2911.3335

# Row 132: On Board for Lar
$ws.Range("H132").Value = 3695.68
$ws.Range("I132").Value = 3466.8096
$ws.Range("K132").Value = 10400.4288
$ws.Range("M132").Value = -7870.4288

$ws = $wb.Worksheets.Item("LTW")
# Row 11: A Thorn in One's Hide
$ws.Range("H11").Value = 2824.5
$ws.Range("J11").Value = 2824.5
$ws.Range("L11").Value = 2824.5
$ws.Range("N11").Value = -3104.5

# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 1775.6111
$ws.Range("I22").Value = 1399.125
$ws.Range("J22").Value = 2076.8
$ws.Range("K22").Value = 1399.125
$ws.Range("L22").Value = 2076.8
$ws.Range("M22").Value = -1104.125
$ws.Range("N22").Value = -2666.8

# Row 25: A Rush on Ringbands
$ws.Range("H25").Value = 2059.3333
$ws.Range("I25").Value = 3119.6667
$ws.Range("J25").Value = 999
$ws.Range("K25").Value = 3119.6667
$ws.Range("L25").Value = 999
$ws.Range("M25").Value = -2889.6667
$ws.Range("N25").Value = -1459

# Row 27: Fire and Hide
$ws.Range("H27").Value = 1775.6111
$ws.Range("I27").Value = 1399.125
$ws.Range("J27").Value = 2076.8
$ws.Range("K27").Value = 1399.125
$ws.Range("L27").Value = 2076.8
$ws.Range("M27").Value = -1292.125
$ws.Range("N27").Value = -2290.8

# Row 110: Breeches of Trust
$ws.Range("H110").Value = 40001
$ws.Range("J110").Value = 40001
$ws.Range("L110").Value = 40001
$ws.Range("N110").Value = -48181

# Row 131: For What Was Gleaned
$ws.Range("H131").Value = 70000
$ws.Range("J131").Value = 70000
$ws.Range("L131").Value = 70000
$ws.Range("N131").Value = -80080

$ws = $wb.Worksheets.Item("WVR")
# Row 3: Trew Enough
$ws.Range("H3").Value = 10000000
$ws.Range("I3").Value = 10000000
$ws.Range("K3").Value = 10000000
$ws.Range("M3").Value = -9999886

# Row 6: Burn Me Up
$ws.Range("H6").Value = 2000
$ws.Range("J6").Value = 2000
$ws.Range("L6").Value = 2000
$ws.Range("N6").Value = -2230
